$wb = $excel.ActiveWorkbook

# Update the "F" column (想去人数 / number of people interested) values on both the
# "展览" and "全部类型" worksheets, which carry duplicate rows.
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 159
    3  = 7105
    4  = 5062
    6  = 159
    7  = 33
    11 = 74
    12 = 190
    13 = 622
    14 = 165
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
